# Auto update Excel log
# Appends new sensor log rows to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 123-135 ---
$wsPIR = $wb.Worksheets.Item("PIR")
$pirData = @(
    @("2026-01-28","18:11:37","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:11:41","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:11:43","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:11:45","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:11:50","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:11:57","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:12:00","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:12:05","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:12:10","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:12:17","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:12:21","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:12:25","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","18:12:30","18:00","Bathroom","No Motion","Inactive")
)
$startRow = 123
for ($i = 0; $i -lt $pirData.Count; $i++) {
    $row = $startRow + $i
    $rec = $pirData[$i]
    # Column A holds a date-looking string ("2026-01-28"); force Text format
    # so Excel stores it verbatim instead of auto-converting to a date serial.
    $wsPIR.Cells.Item($row, 1).NumberFormat = "@"
    for ($c = 1; $c -le 6; $c++) {
        $wsPIR.Cells.Item($row, $c).Value = $rec[$c - 1]
    }
}

# --- Humidity sheet: append rows 119-131 ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityData = @(
    @("2026-01-28","18:11:37","18:00","Bathroom","88.3%","Active"),
    @("2026-01-28","18:11:39","18:00","Bathroom","87.3%","Active"),
    @("2026-01-28","18:11:42","18:00","Bathroom","88.2%","Active"),
    @("2026-01-28","18:11:47","18:00","Bathroom","87.3%","Active"),
    @("2026-01-28","18:11:51","18:00","Bathroom","88.2%","Active"),
    @("2026-01-28","18:11:55","18:00","Bathroom","87.4%","Active"),
    @("2026-01-28","18:12:07","18:00","Bathroom","88.3%","Active"),
    @("2026-01-28","18:12:11","18:00","Bathroom","88.3%","Active"),
    @("2026-01-28","18:12:15","18:00","Bathroom","87.3%","Active"),
    @("2026-01-28","18:12:19","18:00","Bathroom","88.2%","Active"),
    @("2026-01-28","18:12:23","18:00","Bathroom","88.2%","Active"),
    @("2026-01-28","18:12:31","18:00","Bathroom","88.2%","Active"),
    @("2026-01-28","18:12:35","18:00","Bathroom","87.3%","Active")
)
$startRow = 119
for ($i = 0; $i -lt $humidityData.Count; $i++) {
    $row = $startRow + $i
    $rec = $humidityData[$i]
    # Column A (date) and Column E (percentage text) both need Text format
    # to avoid Excel auto-converting them to a date serial / numeric percent.
    $wsHumidity.Cells.Item($row, 1).NumberFormat = "@"
    $wsHumidity.Cells.Item($row, 5).NumberFormat = "@"
    for ($c = 1; $c -le 6; $c++) {
        $wsHumidity.Cells.Item($row, $c).Value = $rec[$c - 1]
    }
}

# --- Temperature sheet: append rows 118-130 ---
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureData = @(
    @("2026-01-28","18:11:36","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:11:38","18:00","Bathroom","23.0C","Active"),
    @("2026-01-28","18:11:40","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:11:42","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:11:48","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:11:52","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:11:56","18:00","Bathroom","23.0C","Active"),
    @("2026-01-28","18:12:08","18:00","Bathroom","23.0C","Active"),
    @("2026-01-28","18:12:12","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:12:16","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:12:20","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:12:24","18:00","Bathroom","22.9C","Active"),
    @("2026-01-28","18:12:32","18:00","Bathroom","22.9C","Active")
)
$startRow = 118
for ($i = 0; $i -lt $temperatureData.Count; $i++) {
    $row = $startRow + $i
    $rec = $temperatureData[$i]
    # Column A holds a date-looking string; force Text format.
    $wsTemperature.Cells.Item($row, 1).NumberFormat = "@"
    for ($c = 1; $c -le 6; $c++) {
        $wsTemperature.Cells.Item($row, $c).Value = $rec[$c - 1]
    }
}
